# Apply updates for "Add data for 2022-08-07"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the "through" date everywhere it's referenced
$ws.Name = "Through 2022-07-30"
$ws.Range("I1").Value = "2022 (through 07-30)"

# Update the July total (row 8) and the grand Total row (row 14) for column I
$ws.Range("I8").Value = 164
$ws.Range("I14").Value = 970
